$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2017-02-15 06:06:13"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("H2").Value = "2017-02-15 06:05:54"
$wsZhCn.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39742656c34e65e9a770f8196cc69173e81f45a0/e2e/30be2b79-a958-481b-bd01-0d152dd008b2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf335ef29c8a88c27250556c7f68f6a7947e25e1/e2e/30be2b79-a958-481b-bd01-0d152dd008b2.md."
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsZhCn.Columns.Item(18).ColumnWidth = 39.16666666666667

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("H2").Value = "2017-02-15 06:06:13"
$wsDeDe.Range("R2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39742656c34e65e9a770f8196cc69173e81f45a0/e2e/30be2b79-a958-481b-bd01-0d152dd008b2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cf335ef29c8a88c27250556c7f68f6a7947e25e1/e2e/30be2b79-a958-481b-bd01-0d152dd008b2.md."
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(18).ColumnWidth = 39.16666666666667
